$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '35.220.12'
$ws.Range("E2").Value = '  +0.18%  '
$ws.Range("D3").Value = '1.893.73'
$ws.Range("E3").Value = '  +2.13%  '
$ws.Range("E4").Value = '  -0.50%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '243.00'
$ws.Range("D5").NumberFormat = "General"
$ws.Range("E5").Value = '  +2.27%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '0.655'
$ws.Range("D6").NumberFormat = "General"
$ws.Range("E6").Value = '  +5.48%  '
$ws.Range("E7").Value = '  -0.44%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '41.35'
$ws.Range("D8").NumberFormat = "General"
$ws.Range("E8").Value = '  -0.84%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.344'
$ws.Range("D9").NumberFormat = "General"
$ws.Range("E9").Value = '  +5.30%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '50.33'
$ws.Range("D10").NumberFormat = "General"
$ws.Range("E10").Value = '  +7.92%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.0707'
$ws.Range("D11").NumberFormat = "General"
$ws.Range("E11").Value = '  +2.26%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.0996'
$ws.Range("D12").NumberFormat = "General"
$ws.Range("E12").Value = '  +0.52%  '
$ws.Range("D13").Value = '2.170.06'
$ws.Range("E13").Value = '  +2.28%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '11.96'
$ws.Range("D14").NumberFormat = "General"
$ws.Range("E14").Value = '  +5.25%  '
$ws.Range("B15").Value = 'Polygon'
$ws.Range("C15").Value = 'https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic'
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '0.691'
$ws.Range("D15").NumberFormat = "General"
$ws.Range("E15").Value = '  +2.69%  '
$ws.Range("B16").Value = 'WrappedEther'
$ws.Range("C16").Value = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
$ws.Range("D16").Value = '1.891.80'
$ws.Range("E16").Value = '  +1.70%  '
$ws.Range("E17").Value = '  +1.80%  '
$ws.Range("D18").Value = '35.213.09'
$ws.Range("E18").Value = '  +0.28%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '71.12'
$ws.Range("D19").NumberFormat = "General"
$ws.Range("E19").Value = '  +1.65%  '
$ws.Range("D20").Value = '0.0₃0811'
$ws.Range("E20").Value = '  +2.67%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '240.57'
$ws.Range("D21").NumberFormat = "General"
$ws.Range("E21").Value = '  +0.12%  '
$ws.Range("E22").Value = '  +2.16%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '4.72'
$ws.Range("D23").NumberFormat = "General"
$ws.Range("E23").Value = '  +0.24%  '
$ws.Range("E24").Value = '  -0.44%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '2.41'
$ws.Range("D25").NumberFormat = "General"
$ws.Range("E25").Value = '  +32.91%  '
$ws.Range("E26").Value = '  +0.30%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '169.71'
$ws.Range("D27").NumberFormat = "General"
$ws.Range("E27").Value = '  +0.57%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '8.34'
$ws.Range("D28").NumberFormat = "General"
$ws.Range("E28").Value = '  +4.67%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '18.11'
$ws.Range("D29").NumberFormat = "General"
$ws.Range("E29").Value = '  +3.31%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '0.126'
$ws.Range("D30").NumberFormat = "General"
$ws.Range("E30").Value = '  +2.26%  '
$ws.Range("E31").Value = '  +3.09%  '
$ws.Range("B32").Value = 'BinanceUSD'
$ws.Range("C32").Value = 'https://coinranking.com/coin/vSo2fu9iE1s0Y+binanceusd-busd'
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '1.01'
$ws.Range("D32").NumberFormat = "General"
$ws.Range("E32").Value = '  -0.44%  '
$ws.Range("B33").Value = 'Hedera'
$ws.Range("C33").Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '0.0558'
$ws.Range("D33").NumberFormat = "General"
$ws.Range("E33").Value = '  +1.21%  '
$ws.Range("B34").Value = 'ImmutableX'
$ws.Range("C34").Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '0.936'
$ws.Range("D34").NumberFormat = "General"
$ws.Range("E34").Value = '  +17.70%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '4.09'
$ws.Range("D35").NumberFormat = "General"
$ws.Range("E35").Value = '  +2.41%  '
$ws.Range("E36").Value = '  -0.70%  '
$ws.Range("E37").Value = '  +1.23%  '
$ws.Range("E38").Value = '  +2.27%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '1.08'
$ws.Range("D39").NumberFormat = "General"
$ws.Range("E39").Value = '  +1.53%  '
$ws.Range("E40").Value = '  +3.11%  '
$ws.Range("E41").Value = '  +14.70%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '15.93'
$ws.Range("D42").NumberFormat = "General"
$ws.Range("E42").Value = '  +8.02%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '88.81'
$ws.Range("D43").NumberFormat = "General"
$ws.Range("E43").Value = '  -0.96%  '
$ws.Range("D44").Value = '1.335.29'
$ws.Range("E44").Value = '  -0.44%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '48.36'
$ws.Range("D45").NumberFormat = "General"
$ws.Range("E45").Value = '  +41.53%  '
$ws.Range("E46").Value = '  +2.77%  '
$ws.Range("E47").Value = '  -1.52%  '
$ws.Range("E48").Value = '  +1.15%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '6.49'
$ws.Range("D49").NumberFormat = "General"
$ws.Range("E49").Value = '  +0.74%  '
$ws.Range("D50").Value = '2.079.03'
$ws.Range("E50").Value = '  +2.43%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '11.21'
$ws.Range("D51").NumberFormat = "General"
$ws.Range("E51").Value = '  -12.86%  '
